$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.633.80"
$ws.Range("E2").Value = "  +3.57%  "

# Row 3
$ws.Range("D3").Value = "1.698.41"
$ws.Range("E3").Value = "  +2.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.90%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3951"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.78%  "

# Row 8
$ws.Range("E8").Value = "  +2.41%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.542"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.82"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08814"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.298"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +8.85%  "

# Row 14
$ws.Range("E14").Value = "  +3.55%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001332"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.654"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.47%  "

# Row 17
$ws.Range("D17").Value = "1.698.33"
$ws.Range("E17").Value = "  +1.83%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07097"
$ws.Range("D19").ClearFormats()

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.907"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.16%  "

# Row 22
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.16"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.07%  "

# Row 24
$ws.Range("D24").Value = "24.622.87"
$ws.Range("E24").Value = "  +3.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.134"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +14.29%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.335"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.07%  "

# Row 27
$ws.Range("E27").Value = "  +3.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.97"
$ws.Range("D28").ClearFormats()

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.247"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.29"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.64%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.638"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +16.56%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.111"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.56%  "

# Row 33
$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.885.98"
$ws.Range("E33").Value = "  +1.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.452"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.67%  "

# Row 35
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.36"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +10.29%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2759"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.88%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.946"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.50%  "

# Row 39
$ws.Range("E39").Value = "  +3.56%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02801"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.36%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09123"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7771"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.39%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.467"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7297"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.31%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.521"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.217"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.385"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +20.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.17%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08044"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.21%  "
